$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per latest crypto data refresh
$ws.Range("D2").Value = "'42.951.36"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "'2.279.41"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'250.39"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("E6").Value = "  +1.15%  "
$ws.Range("D7").Value = "'77.81"
$ws.Range("E7").Value = "  +8.36%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "'0.653"
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("D10").Value = "'40.66"
$ws.Range("E10").Value = "  +5.44%  "
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").Value = "'7.33"
$ws.Range("E12").Value = "  -0.62%  "
$ws.Range("D13").Value = "'0.106"
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("D14").Value = "'2.618.06"
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("D15").Value = "'15.11"
$ws.Range("E15").Value = "  +1.36%  "
$ws.Range("D16").Value = "'0.871"
$ws.Range("E16").Value = "  -1.70%  "
$ws.Range("D17").Value = "'2.285.69"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("D18").Value = "'42.832.39"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").Value = "'0.0₃0995"
$ws.Range("E19").Value = "  -1.87%  "
$ws.Range("D20").Value = "'6.23"
$ws.Range("E20").Value = "  -1.22%  "
$ws.Range("D21").Value = "'72.15"
$ws.Range("E21").Value = "  -1.79%  "
$ws.Range("D22").Value = "'234.01"
$ws.Range("E22").Value = "  -0.94%  "
$ws.Range("D23").Value = "'2.15"
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("D24").Value = "'3.79"
$ws.Range("E24").Value = "  -5.45%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").Value = "'11.37"
$ws.Range("E26").Value = "  -1.30%  "
$ws.Range("D27").Value = "'2.36"
$ws.Range("E27").Value = "  -3.00%  "
$ws.Range("E28").Value = "  -1.35%  "
$ws.Range("D29").Value = "'168.25"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").Value = "'20.94"
$ws.Range("E30").Value = "  -0.21%  "
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("D32").Value = "'0.0852"
$ws.Range("E32").Value = "  +6.81%  "
$ws.Range("E33").Value = "  -4.16%  "
$ws.Range("D34").Value = "'30.46"
$ws.Range("E34").Value = "  -2.92%  "
$ws.Range("E35").Value = "  +0.99%  "
$ws.Range("D36").Value = "'4.59"
$ws.Range("E36").Value = "  +0.81%  "
$ws.Range("D37").Value = "'4.74"
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("D38").Value = "'0.0305"
$ws.Range("E38").Value = "  -2.15%  "
$ws.Range("D39").Value = "'13.71"
$ws.Range("E39").Value = "  +3.27%  "
$ws.Range("D40").Value = "'2.27"
$ws.Range("E40").Value = "  -2.64%  "
$ws.Range("D41").Value = "'5.86"
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("D42").Value = "'112.80"
$ws.Range("E42").Value = "  +18.13%  "
$ws.Range("D43").Value = "'0.209"
$ws.Range("E43").Value = "  -0.60%  "
$ws.Range("D44").Value = "'61.23"
$ws.Range("E44").Value = "  -1.09%  "
$ws.Range("D45").Value = "'8.88"
$ws.Range("E45").Value = "  -3.13%  "
$ws.Range("E46").Value = "  -1.90%  "
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("D48").Value = "'4.54"
$ws.Range("E48").Value = "  -9.53%  "
$ws.Range("D49").Value = "'1.16"
$ws.Range("E49").Value = "  -2.62%  "
$ws.Range("D50").Value = "'1.17"
$ws.Range("E50").Value = "  -2.31%  "
$ws.Range("D51").Value = "'4.25"
$ws.Range("E51").Value = "  +0.28%  "
